# Actualización de tabla champion: se agrega una columna "component" al
# inicio de la tabla, desplazando "model"/"train_rmse"/"test_rmse" una
# columna a la derecha, y se completa con los nuevos valores de cada
# componente del PIB.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insertar una nueva columna antes de la columna A, desplazando el resto
# de las columnas (B,C) hacia la derecha (C,D).
$ws.Range("A1").EntireColumn.Insert()

# Cabecera + datos de la nueva columna "component".
$ws.Range("A1").Value = "component"
$components = @("consumo_privado", "exportaciones", "gasto_publico", "importaciones", "inversiones", "variacion")
for ($i = 0; $i -lt $components.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $components[$i]
}

# Formato de la nueva columna: negrita, bordes finos y texto centrado,
# igual que el resto de encabezados de la tabla.
$rng = $ws.Range("A1:A7")
$rng.Font.Bold = $true
$rng.Borders.LineStyle = 1
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160

$wb.Save()
